$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (Muhammad Luqman_20251202_121800) - removed via admin panel,
# all subsequent rows shift up by one.
$ws.Rows.Item(11).Delete()
